# Update the cryptocurrency "Price" (column D) and "Volume(1h)" (column E) values
# to the latest scraped figures, row by row (row 2 = Bitcoin ... row 51 = last coin).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.973.17"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").Value = "2.927.34"
$ws.Range("E3").Value = "  +1.29%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'358.56"
$ws.Range("E5").Value = "  +2.21%  "
$ws.Range("D6").Value = "'110.33"
$ws.Range("E6").Value = "  -0.90%  "
$ws.Range("E7").Value = "  +1.64%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +1.20%  "
$ws.Range("D10").Value = "'39.51"
$ws.Range("E10").Value = "  -0.64%  "
$ws.Range("D11").Value = "'0.0883"
$ws.Range("E11").Value = "  +3.02%  "
$ws.Range("E12").Value = "  +0.56%  "
$ws.Range("E13").Value = "  -1.10%  "
$ws.Range("D14").Value = "'7.89"
$ws.Range("E14").Value = "  +1.70%  "
$ws.Range("D15").Value = "3.390.71"
$ws.Range("E15").Value = "  +1.39%  "
$ws.Range("D16").Value = "2.925.74"
$ws.Range("E16").Value = "  +1.34%  "
$ws.Range("D17").Value = "'0.987"
$ws.Range("E17").Value = "  -0.84%  "
$ws.Range("D18").Value = "51.988.19"
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("D19").Value = "'3.32"
$ws.Range("E19").Value = "  -0.26%  "
$ws.Range("D20").Value = "'7.59"
$ws.Range("E20").Value = "  -1.28%  "
$ws.Range("D21").Value = "'14.07"
$ws.Range("E21").Value = "  -2.47%  "
$ws.Range("E22").Value = "  +0.44%  "
$ws.Range("D23").Value = "'71.00"
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("D24").Value = "'270.69"
$ws.Range("E24").Value = "  +0.53%  "
$ws.Range("D25").Value = "'2.83"
$ws.Range("E25").Value = "  +2.12%  "
$ws.Range("D26").Value = "'0.184"
$ws.Range("E26").Value = "  +12.08%  "
$ws.Range("D27").Value = "'27.19"
$ws.Range("E27").Value = "  +2.64%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("D29").Value = "'7.42"
$ws.Range("E29").Value = "  +15.38%  "
$ws.Range("D30").Value = "'0.108"
$ws.Range("E30").Value = "  +14.52%  "
$ws.Range("E31").Value = "  +1.08%  "
$ws.Range("D32").Value = "'38.64"
$ws.Range("E32").Value = "  +1.01%  "
$ws.Range("D33").Value = "'6.04"
$ws.Range("E33").Value = "  -1.88%  "
$ws.Range("D34").Value = "'52.24"
$ws.Range("E34").Value = "  -1.09%  "
$ws.Range("E35").Value = "  -2.26%  "
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("E37").Value = "  -13.87%  "
$ws.Range("D38").Value = "'3.25"
$ws.Range("E38").Value = "  -1.30%  "
$ws.Range("E39").Value = "  -0.15%  "
$ws.Range("E40").Value = "  -1.01%  "
$ws.Range("D41").Value = "'2.75"
$ws.Range("E41").Value = "  +4.31%  "
$ws.Range("E42").Value = "  +3.13%  "
$ws.Range("D43").Value = "'23.25"
$ws.Range("E43").Value = "  +2.40%  "
$ws.Range("D44").Value = "'119.29"
$ws.Range("E44").Value = "  -2.47%  "
$ws.Range("E45").Value = "  -1.51%  "
$ws.Range("E46").Value = "  -0.09%  "
$ws.Range("E47").Value = "  -2.45%  "
$ws.Range("D48").Value = "2.139.79"
$ws.Range("E48").Value = "  -2.49%  "
$ws.Range("E49").Value = "  -8.02%  "
$ws.Range("E50").Value = "  +2.86%  "
$ws.Range("D51").Value = "'9.21"
$ws.Range("E51").Value = "  +1.50%  "
